$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 9
$ws.Cells.Item($row, 1).Value = 1746359037
$ws.Cells.Item($row, 2).Value = "update"
$ws.Cells.Item($row, 3).Value = "variable"
$ws.Cells.Item($row, 4).Value = "dep_sante___variable_3"
$ws.Cells.Item($row, 6).Value = "type"
$ws.Cells.Item($row, 7).Value = "integer"
$ws.Cells.Item($row, 8).Value = "string"
